# Regenerate sval data to filter save games.
# Update columns B, C, D, E (and derived sum column G) for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    3 = @{ B = 0.04271373187048222; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 2.342150324931327 }
    4 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 }
    5 = @{ B = 0.1190320826869504; C = 1.655778082260271; D = 0.1494219747398047; E = 0.4942365360607697; G = 2.418468675747795 }
    6 = @{ B = 3.286832544864788;  C = 1.655778082260271; D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    7 = @{ B = 1.455362044514542;  C = 1.655778082260271; D = 0.7527432677738641; E = 10.19245300693656;  G = 14.05633640148523 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
